$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.83880576904272
$ws.Range("C2").Value = 10.62394913076492
$ws.Range("D2").Value = 4.818927982790526
$ws.Range("F2").Value = 29.30690430826319
$ws.Range("G2").Value = 3.626321559908817
$ws.Range("I2").Value = 22.1744173994536
$ws.Range("M2").Value = 19.96912112489786
$ws.Range("N2").Value = 17.56673258165777

$ws.Range("B3").Value = 13.27768040354618
$ws.Range("C3").Value = 10.01559014073007
$ws.Range("D3").Value = 4.842248384897366
$ws.Range("F3").Value = 28.94125244399019
$ws.Range("G3").Value = 3.630003277128562
$ws.Range("I3").Value = 22.12318453253697
$ws.Range("M3").Value = 19.3658123765279
$ws.Range("N3").Value = 17.63548376069562

$ws.Range("B4").Value = 12.92464632814239
$ws.Range("C4").Value = 9.625036476472584
$ws.Range("D4").Value = 4.857157732237556
$ws.Range("F4").Value = 28.72489357479759
$ws.Range("G4").Value = 3.632379305103977
$ws.Range("I4").Value = 22.0981593658074
$ws.Range("M4").Value = 18.99400431950565
$ws.Range("N4").Value = 17.6796392945844

$ws.Range("B5").Value = 12.77888903744139
$ws.Range("C5").Value = 9.461765965658399
$ws.Range("D5").Value = 4.86338232631677
$ws.Range("F5").Value = 28.63887693845233
$ws.Range("G5").Value = 3.633376700472912
$ws.Range("I5").Value = 22.08957876943155
$ws.Range("M5").Value = 18.84239795604602
$ws.Range("N5").Value = 17.6981221897242

$ws.Range("B6").Value = 12.75457949112997
$ws.Range("C6").Value = 9.434411794355306
$ws.Range("D6").Value = 4.864424921922851
$ws.Range("F6").Value = 28.62472672901767
$ws.Range("G6").Value = 3.633544080902466
$ws.Range("I6").Value = 22.08825161614434
$ws.Range("M6").Value = 18.81722596690726
$ws.Range("N6").Value = 17.70122082852948

$ws.Range("B7").Value = 12.92268791620947
$ws.Range("C7").Value = 9.622850982373638
$ws.Range("D7").Value = 4.85724107577281
$ws.Range("F7").Value = 28.72372468379453
$ws.Range("G7").Value = 3.632392638165346
$ws.Range("I7").Value = 22.09803709757834
$ws.Range("M7").Value = 18.99195972906617
$ws.Range("N7").Value = 17.6798865796762

$ws.Range("B8").Value = 13.64723255319886
$ws.Range("C8").Value = 10.41780709249944
$ws.Range("D8").Value = 4.826846547780032
$ws.Range("F8").Value = 29.17919778736898
$ws.Range("G8").Value = 3.627567127121262
$ws.Range("I8").Value = 22.15541786605014
$ws.Range("M8").Value = 19.76153378576642
$ws.Range("N8").Value = 17.59003541198721

$ws.Range("B9").Value = 14.99120550088236
$ws.Range("C9").Value = 11.83602944588243
$ws.Range("D9").Value = 4.771908894810785
$ws.Range("F9").Value = 30.13237725487124
$ws.Range("G9").Value = 3.619014938999003
$ws.Range("I9").Value = 22.31889814294994
$ws.Range("M9").Value = 21.24931167273584
$ws.Range("N9").Value = 17.42920816820067

$ws.Range("B10").Value = 15.92093389650063
$ws.Range("C10").Value = 12.78638706213128
$ws.Range("D10").Value = 4.734362662053145
$ws.Range("F10").Value = 30.86271384207194
$ws.Range("G10").Value = 3.613279315399895
$ws.Range("I10").Value = 22.46985213146387
$ws.Range("M10").Value = 22.31649042186726
$ws.Range("N10").Value = 17.32036315213719

$ws.Range("B11").Value = 16.32953678030002
$ws.Range("C11").Value = 13.19799609199176
$ws.Range("D11").Value = 4.717887626976293
$ws.Range("F11").Value = 31.20005010024196
$ws.Range("G11").Value = 3.610787351321425
$ws.Range("I11").Value = 22.5451378789012
$ws.Range("M11").Value = 22.7939763778209
$ws.Range("N11").Value = 17.27285790264667

$ws.Range("B12").Value = 16.4820775824043
$ws.Range("C12").Value = 13.35083237882625
$ws.Range("D12").Value = 4.711735543885083
$ws.Range("F12").Value = 31.32840344908172
$ws.Range("G12").Value = 3.609860438851339
$ws.Range("I12").Value = 22.57458787285457
$ws.Range("M12").Value = 22.97346615961055
$ws.Range("N12").Value = 17.25515704627144

$ws.Range("B13").Value = 16.44932434351489
$ws.Range("C13").Value = 13.31805187069756
$ws.Range("D13").Value = 4.713056656202998
$ws.Range("F13").Value = 31.30073500285033
$ws.Range("G13").Value = 3.610059323305887
$ws.Range("I13").Value = 22.56820362719768
$ws.Range("M13").Value = 22.93487132160792
$ws.Range("N13").Value = 17.25895642587633

$ws.Range("B14").Value = 16.34213092226253
$ws.Range("C14").Value = 13.21063103213679
$ws.Range("D14").Value = 4.717379756567526
$ws.Range("F14").Value = 31.21059825091365
$ws.Range("G14").Value = 3.610710758821626
$ws.Range("I14").Value = 22.54754195637539
$ws.Range("M14").Value = 22.80877058882399
$ws.Range("N14").Value = 17.27139586645533

$ws.Range("B15").Value = 16.27618342600602
$ws.Range("C15").Value = 13.14443657052251
$ws.Range("D15").Value = 4.720039055736545
$ws.Range("F15").Value = 31.15546289860083
$ws.Range("G15").Value = 3.611111958609165
$ws.Range("I15").Value = 22.5350082812168
$ws.Range("M15").Value = 22.73135291747701
$ws.Range("N15").Value = 17.27905291974363

$ws.Range("B16").Value = 15.89393101444134
$ws.Range("C16").Value = 12.75906586440137
$ws.Range("D16").Value = 4.735451488342661
$ws.Range("F16").Value = 30.84076032810996
$ws.Range("G16").Value = 3.613444519814177
$ws.Range("I16").Value = 22.46506430777695
$ws.Range("M16").Value = 22.28510955026594
$ws.Range("N16").Value = 17.32350811020474

$ws.Range("B17").Value = 15.65566382018756
$ws.Range("C17").Value = 12.51730685131285
$ws.Range("D17").Value = 4.745061213182844
$ws.Range("F17").Value = 30.64892015164488
$ws.Range("G17").Value = 3.614905407831107
$ws.Range("I17").Value = 22.42384358183664
$ws.Range("M17").Value = 22.00917894959322
$ws.Range("N17").Value = 17.35129406984759

$ws.Range("B18").Value = 15.51727588353591
$ws.Range("C18").Value = 12.37630609297949
$ws.Range("D18").Value = 4.750645431225462
$ws.Range("F18").Value = 30.53906633325701
$ws.Range("G18").Value = 3.615756710129826
$ws.Range("I18").Value = 22.40075798140294
$ws.Range("M18").Value = 21.849729359504
$ws.Range("N18").Value = 17.36746493485202

$ws.Range("B19").Value = 15.47019361974868
$ws.Range("C19").Value = 12.32823300078462
$ws.Range("D19").Value = 4.75254594453561
$ws.Range("F19").Value = 30.50195919695125
$ws.Range("G19").Value = 3.616046845716796
$ws.Range("I19").Value = 22.39304896538744
$ws.Range("M19").Value = 21.79562076126028
$ws.Range("N19").Value = 17.37297260927826

$ws.Range("B20").Value = 15.68116773647737
$ws.Range("C20").Value = 12.54324447540452
$ws.Range("D20").Value = 4.744032349044558
$ws.Range("F20").Value = 30.66929229192202
$ws.Range("G20").Value = 3.614748752219628
$ws.Range("I20").Value = 22.42816714186359
$ws.Range("M20").Value = 22.03863032049294
$ws.Range("N20").Value = 17.34831663816497

$ws.Range("B21").Value = 16.37367653863063
$ws.Range("C21").Value = 13.2422657598466
$ws.Range("D21").Value = 4.716107607992612
$ws.Range("F21").Value = 31.23705797613895
$ws.Range("G21").Value = 3.610518962993512
$ws.Range("I21").Value = 22.55358534572235
$ws.Range("M21").Value = 22.84584667673693
$ws.Range("N21").Value = 17.26773427947002

$ws.Range("B22").Value = 16.81346977207273
$ws.Range("C22").Value = 13.68143249500827
$ws.Range("D22").Value = 4.698362122479534
$ws.Range("F22").Value = 31.61163364476741
$ws.Range("G22").Value = 3.607852073860698
$ws.Range("I22").Value = 22.64103170596056
$ws.Range("M22").Value = 23.36562366271158
$ws.Range("N22").Value = 17.2167495598014

$ws.Range("B23").Value = 16.57995251623872
$ws.Range("C23").Value = 13.44867341397506
$ws.Range("D23").Value = 4.707787126164313
$ws.Range("F23").Value = 31.4114347112501
$ws.Range("G23").Value = 3.60926655673076
$ws.Range("I23").Value = 22.59386257183389
$ws.Range("M23").Value = 23.08897504387128
$ws.Range("N23").Value = 17.24380746158104

$ws.Range("B24").Value = 15.66964178222332
$ws.Range("C24").Value = 12.5315243330224
$ws.Range("D24").Value = 4.744497313293942
$ws.Range("F24").Value = 30.66008067482085
$ws.Range("G24").Value = 3.614819540648079
$ws.Range("I24").Value = 22.42621055104233
$ws.Range("M24").Value = 22.02531787572816
$ws.Range("N24").Value = 17.34966212341259

$ws.Range("B25").Value = 14.63706425540825
$ws.Range("C25").Value = 11.46810963326184
$ws.Range("D25").Value = 4.786274402616248
$ws.Range("F25").Value = 29.86875864427717
$ws.Range("G25").Value = 3.621231817119506
$ws.Range("I25").Value = 22.26923246303584
$ws.Range("M25").Value = 20.85050875550169
$ws.Range("N25").Value = 17.47107653257314
